$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPME Acetone Water")

# Fill in the raw measured values for row 12 (D12:I12)
$ws.Range("D12").Value = 0.15581490000000001
$ws.Range("E12").Value = 0.38922275000000001
$ws.Range("F12").Value = 0.45496236000000001
$ws.Range("G12").Value = 0.019652340000000001
$ws.Range("H12").Value = 0.20190653
$ws.Range("I12").Value = 0.77844111999999999

# Fill in the computed-ratio formulas for row 12 (J12:U12), matching the
# pattern already used in row 11 just above.
$ws.Range("P12").Formula = '=(D12*$A$4)/((D12*$A$4)+(E12*$B$4)+(F12*$C$4))'
$ws.Range("Q12").Formula = '=(E12*$B$4)/((D12*$A$4)+(E12*$B$4)+(F12*$C$4))'
$ws.Range("R12").Formula = '=(F12*$C$4)/((D12*$A$4)+(E12*$B$4)+(F12*$C$4))'
$ws.Range("S12").Formula = '=(G12*$A$4)/((G12*$A$4)+(H12*$B$4)+(I12*$C$4))'
$ws.Range("T12").Formula = '=(H12*$B$4)/((G12*$A$4)+(H12*$B$4)+(I12*$C$4))'
$ws.Range("U12").Formula = '=(I12*$C$4)/((G12*$A$4)+(H12*$B$4)+(I12*$C$4))'

$ws.Range("J12").Formula = '=(P12/$A$6)/((P12/$A$6)+(Q12/$B$6)+(R12/$C$6))'
$ws.Range("K12").Formula = '=(Q12/$B$6)/((P12/$A$6)+(Q12/$B$6)+(R12/$C$6))'
$ws.Range("L12").Formula = '=(R12/$C$6)/((P12/$A$6)+(Q12/$B$6)+(R12/$C$6))'
$ws.Range("M12").Formula = '=(S12/$A$6)/((S12/$A$6)+(T12/$B$6)+(U12/$C$6))'
$ws.Range("N12").Formula = '=(T12/$B$6)/((S12/$A$6)+(T12/$B$6)+(U12/$C$6))'
$ws.Range("O12").Formula = '=(U12/$C$6)/((S12/$A$6)+(T12/$B$6)+(U12/$C$6))'

# Recalculate so the cached formula results match the newly entered values
$excel.Calculate()

# Update the selected cell on this sheet to reflect the author's final
# cursor position after the edit.
$ws.Activate()
$null = $ws.Range("I13").Select()
